$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a numeric-looking string must be forced to Text format
# so Excel stores them as text (matching the source inline-string cells) instead of
# silently parsing them into floating point numbers.

$ws.Range("D2").Value = "26.669.37"
$ws.Range("E2").Value = "  -1.55%  "
$ws.Range("D3").Value = "1.594.48"
$ws.Range("E3").Value = "  -1.77%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.40"
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0618"
$ws.Range("E8").Value = "  -1.87%  "
$ws.Range("E9").Value = "  -1.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.66"
$ws.Range("E10").Value = "  -1.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0837"
$ws.Range("E11").Value = "  -0.98%  "
$ws.Range("D12").Value = "1.817.87"
$ws.Range("E12").Value = "  -1.78%  "
$ws.Range("D13").Value = "1.615.38"
$ws.Range("E13").Value = "  -0.67%  "
$ws.Range("E14").Value = "  -2.63%  "
$ws.Range("E15").Value = "  -3.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.75"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").Value = "26.637.52"
$ws.Range("E18").Value = "  -1.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "209.26"
$ws.Range("E19").Value = "  -1.88%  "
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("E21").Value = "  -2.82%  "
$ws.Range("E22").Value = "  -2.52%  "
$ws.Range("E23").Value = "  -1.62%  "
$ws.Range("E24").Value = "  -2.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.65"
$ws.Range("E25").Value = "  -1.06%  "
$ws.Range("E27").Value = "  -3.13%  "
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("E29").Value = "  -1.35%  "
$ws.Range("E30").Value = "  -1.34%  "
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.22"
$ws.Range("E32").Value = "  -3.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.662"
$ws.Range("E33").Value = "  -8.91%  "
$ws.Range("E34").Value = "  -2.64%  "
$ws.Range("D35").Value = "1.285.37"
$ws.Range("E35").Value = "  -5.53%  "
$ws.Range("E36").Value = "  -0.88%  "
$ws.Range("E37").Value = "  -5.65%  "
$ws.Range("E38").Value = "  -3.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.834"
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("E41").Value = "  -1.08%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.36"
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.20"
$ws.Range("E43").Value = "  -1.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.47"
$ws.Range("E44").Value = "  -1.37%  "
$ws.Range("D45").Value = "1.730.61"
$ws.Range("E45").Value = "  -1.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.71"
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.872"
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("E48").Value = "  -1.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0982"
$ws.Range("E49").Value = "  -2.69%  "
$ws.Range("E50").Value = "  -1.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.53"
$ws.Range("E51").Value = "  -1.68%  "
